# Auto-generated edit script applying the diff to 合肥-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

function Set-CellText {
    param($ws, $row, $col, $val)
    # Force text (no date/number auto-coercion) using Excel's leading-apostrophe convention,
    # then strip the apostrophe back off so the stored value is the literal text itself.
    $ws.Cells.Item($row, $col).Value2 = "'" + $val
}

function Set-CellValue {
    param($ws, $row, $col, $val)
    $ws.Cells.Item($row, $col).Value2 = $val
}

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
# Row 13 is new -> copy column-A style (bold/border/center) from the row above
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
Set-CellValue $ws 2 6 703
Set-CellValue $ws 3 6 31
Set-CellValue $ws 4 6 535
Set-CellValue $ws 9 6 4448
Set-CellValue $ws 10 6 4336
Set-CellText $ws 11 2 '2024-10-04'
Set-CellValue $ws 11 3 '合肥·乐帮•崩原铁绝only同人首展'
Set-CellValue $ws 11 4 '丹霞路488号金星商业城三楼 迷鹿轰趴'
Set-CellValue $ws 11 5 '2024.10.04 10:00-10.05 16:30'
Set-CellValue $ws 11 6 4
Set-CellValue $ws 11 7 58
Set-CellValue $ws 11 8 'https://show.bilibili.com/platform/detail.html?id=91524'
Set-CellValue $ws 11 9 '//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'
Set-CellText $ws 12 2 '2024-10-06'
Set-CellValue $ws 12 3 '合肥·首届火影忍者同人only'
Set-CellValue $ws 12 4 '长江东路金太阳家具广场南门二楼 优极篮球馆'
Set-CellValue $ws 12 5 '2024.10.06 09:30-10.06 17:30'
Set-CellValue $ws 12 6 13
Set-CellValue $ws 12 7 75
Set-CellValue $ws 12 8 'https://show.bilibili.com/platform/detail.html?id=91658'
Set-CellValue $ws 12 9 '//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'
Set-CellValue $ws 13 1 12
Set-CellText $ws 13 2 '2024-10-26'
Set-CellValue $ws 13 3 '合肥·W·A第五人格同人only2.0'
Set-CellValue $ws 13 4 '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
Set-CellValue $ws 13 5 '2024.10.26 09:30-10.26 17:00'
Set-CellValue $ws 13 6 133
Set-CellValue $ws 13 7 68
Set-CellValue $ws 13 8 'https://show.bilibili.com/platform/detail.html?id=91123'
Set-CellValue $ws 13 9 '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
Set-CellValue $ws 2 6 66

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
# Row 16 is new -> copy column-A style (bold/border/center) from the row above
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
Set-CellValue $ws 2 6 703
Set-CellValue $ws 3 6 31
Set-CellValue $ws 4 6 535
Set-CellValue $ws 9 6 4448
Set-CellValue $ws 10 6 4336
Set-CellText $ws 11 2 '2024-10-04'
Set-CellValue $ws 11 3 '合肥·乐帮•崩原铁绝only同人首展'
Set-CellValue $ws 11 4 '丹霞路488号金星商业城三楼 迷鹿轰趴'
Set-CellValue $ws 11 5 '2024.10.04 10:00-10.05 16:30'
Set-CellValue $ws 11 6 4
Set-CellValue $ws 11 7 58
Set-CellValue $ws 11 8 'https://show.bilibili.com/platform/detail.html?id=91524'
Set-CellValue $ws 11 9 '//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'
Set-CellText $ws 12 2 '2024-10-06'
Set-CellValue $ws 12 3 '合肥·首届火影忍者同人only'
Set-CellValue $ws 12 4 '长江东路金太阳家具广场南门二楼 优极篮球馆'
Set-CellValue $ws 12 5 '2024.10.06 09:30-10.06 17:30'
Set-CellValue $ws 12 6 13
Set-CellValue $ws 12 7 75
Set-CellValue $ws 12 8 'https://show.bilibili.com/platform/detail.html?id=91658'
Set-CellValue $ws 12 9 '//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'
Set-CellValue $ws 13 3 '合肥·W·A第五人格同人only2.0'
Set-CellValue $ws 13 4 '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
Set-CellValue $ws 13 5 '2024.10.26 09:30-10.26 17:00'
Set-CellValue $ws 13 6 133
Set-CellValue $ws 13 7 68
Set-CellValue $ws 13 8 'https://show.bilibili.com/platform/detail.html?id=91123'
Set-CellValue $ws 13 9 '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'
Set-CellText $ws 14 2 '2024-10-26'
Set-CellValue $ws 14 3 '合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集'
Set-CellValue $ws 14 5 '2024.10.26 19:30-10.26 21:00'
Set-CellValue $ws 14 6 66
Set-CellValue $ws 14 7 80
Set-CellValue $ws 14 8 'https://show.bilibili.com/platform/detail.html?id=90322'
Set-CellValue $ws 14 9 '//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg'
Set-CellText $ws 15 2 '2024-11-09'
Set-CellValue $ws 15 3 '合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会'
Set-CellValue $ws 15 5 '2024.11.09 19:30-11.09 21:00'
Set-CellValue $ws 15 6 5
Set-CellValue $ws 15 7 64
Set-CellValue $ws 15 8 'https://show.bilibili.com/platform/detail.html?id=90593'
Set-CellValue $ws 15 9 '//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg'
Set-CellValue $ws 16 1 15
Set-CellText $ws 16 2 '2024-12-07'
Set-CellValue $ws 16 3 '合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会'
Set-CellValue $ws 16 4 '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
Set-CellValue $ws 16 5 '2024.12.07 19:30-12.07 21:00'
Set-CellValue $ws 16 6 0
Set-CellValue $ws 16 7 56
Set-CellValue $ws 16 8 'https://show.bilibili.com/platform/detail.html?id=91608'
Set-CellValue $ws 16 9 '//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg'

